$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("L2").Value = 4917
$ws.Range("L3").Value = 5294
$ws.Range("H4").Value = 1764
$ws.Range("L4").Value = 1294
$ws.Range("L5").Value = 312
$ws.Range("L6").Value = 4469
$ws.Range("H7").Value = 26080
$ws.Range("L7").Value = 16286

$ws = $wb.Worksheets.Item(2)
$ws.Range("L7").Value = 531
$ws.Range("L8").Value = 1084
$ws.Range("L11").Value = 264
$ws.Range("L15").Value = 121
$ws.Range("L18").Value = 116
$ws.Range("L19").Value = 446
$ws.Range("L21").Value = 51
$ws.Range("L27").Value = 147
$ws.Range("L29").Value = 891
$ws.Range("L33").Value = 748
$ws.Range("L36").Value = 213
$ws.Range("L37").Value = 616
$ws.Range("L42").Value = 531
$ws.Range("L49").Value = 84
$ws.Range("L51").Value = 207
$ws.Range("L54").Value = 344
$ws.Range("L55").Value = 156
$ws.Range("L59").Value = 30
$ws.Range("L60").Value = 105
$ws.Range("H63").Value = 315
$ws.Range("L63").Value = 46
$ws.Range("L65").Value = 320
$ws.Range("L66").Value = 42
$ws.Range("L67").Value = 566
$ws.Range("L68").Value = 52
$ws.Range("L76").Value = 250
$ws.Range("L78").Value = 210
$ws.Range("L79").Value = 429
$ws.Range("L85").Value = 832
$ws.Range("L90").Value = 163
$ws.Range("L91").Value = 222
$ws.Range("L92").Value = 49
$ws.Range("L94").Value = 200
$ws.Range("L96").Value = 184
$ws.Range("L99").Value = 281
$ws.Range("H101").Value = 26080
$ws.Range("L101").Value = 16286

$ws = $wb.Worksheets.Item(4)
$ws.Range("L3").Value = 55
$ws.Range("L7").Value = 184

$ws = $wb.Worksheets.Item(5)
$ws.Range("L3").Value = 176
$ws.Range("L7").Value = 531

$ws = $wb.Worksheets.Item(6)
$ws.Range("L6").Value = 61
$ws.Range("L7").Value = 264

$ws = $wb.Worksheets.Item(8)
$ws.Range("L3").Value = 340
$ws.Range("L7").Value = 832

$ws = $wb.Worksheets.Item(12)
$ws.Range("L2").Value = 316
$ws.Range("L5").Value = 39
$ws.Range("L6").Value = 280
$ws.Range("L7").Value = 1084

$ws = $wb.Worksheets.Item(14)
$ws.Range("L2").Value = 203
$ws.Range("L6").Value = 225
$ws.Range("L7").Value = 748

$ws = $wb.Worksheets.Item(16)
$ws.Range("L2").Value = 185
$ws.Range("L4").Value = 35
$ws.Range("L6").Value = 171
$ws.Range("L7").Value = 616

$ws = $wb.Worksheets.Item(17)
$ws.Range("L2").Value = 116
$ws.Range("L7").Value = 320

$ws = $wb.Worksheets.Item(18)
$ws.Range("L2").Value = 76
$ws.Range("L6").Value = 58
$ws.Range("L7").Value = 281

$ws = $wb.Worksheets.Item(21)
$ws.Range("L2").Value = 166
$ws.Range("L3").Value = 217
$ws.Range("L6").Value = 131
$ws.Range("L7").Value = 566

$ws = $wb.Worksheets.Item(23)
$ws.Range("L2").Value = 29
$ws.Range("L3").Value = 12
$ws.Range("L7").Value = 84

$ws = $wb.Worksheets.Item(24)
$ws.Range("L6").Value = 168
$ws.Range("L7").Value = 344

$ws = $wb.Worksheets.Item(25)
$ws.Range("L2").Value = 268
$ws.Range("L3").Value = 337
$ws.Range("L4").Value = 44
$ws.Range("L7").Value = 891

$ws = $wb.Worksheets.Item(27)
$ws.Range("L6").Value = 126
$ws.Range("L7").Value = 446

$ws = $wb.Worksheets.Item(29)
$ws.Range("L4").Value = 32
$ws.Range("L7").Value = 250

$ws = $wb.Worksheets.Item(32)
$ws.Range("L2").Value = 153
$ws.Range("L3").Value = 176
$ws.Range("L6").Value = 148
$ws.Range("L7").Value = 531

$ws = $wb.Worksheets.Item(35)
$ws.Range("L6").Value = 64
$ws.Range("L7").Value = 210

$ws = $wb.Worksheets.Item(36)
$ws.Range("L3").Value = 50
$ws.Range("L7").Value = 156

$ws = $wb.Worksheets.Item(40)
$ws.Range("L2").Value = 79
$ws.Range("L7").Value = 222

$ws = $wb.Worksheets.Item(41)
$ws.Range("L6").Value = 27
$ws.Range("L7").Value = 51

$ws = $wb.Worksheets.Item(42)
$ws.Range("L2").Value = 141
$ws.Range("L7").Value = 429

$ws = $wb.Worksheets.Item(45)
$ws.Range("L3").Value = 41
$ws.Range("L7").Value = 116

$ws = $wb.Worksheets.Item(47)
$ws.Range("L3").Value = 65
$ws.Range("L7").Value = 213

$ws = $wb.Worksheets.Item(51)
$ws.Range("L3").Value = 47
$ws.Range("L7").Value = 200

$ws = $wb.Worksheets.Item(54)
$ws.Range("L3").Value = 41
$ws.Range("L7").Value = 121

$ws = $wb.Worksheets.Item(59)
$ws.Range("L6").Value = 14
$ws.Range("L7").Value = 42

$ws = $wb.Worksheets.Item(63)
$ws.Range("L4").Value = 2
$ws.Range("L7").Value = 30

$ws = $wb.Worksheets.Item(66)
$ws.Range("L2").Value = 19
$ws.Range("L6").Value = 19
$ws.Range("L7").Value = 49

$ws = $wb.Worksheets.Item(71)
$ws.Range("L6").Value = 43
$ws.Range("L7").Value = 147

$ws = $wb.Worksheets.Item(74)
$ws.Range("L3").Value = 49
$ws.Range("L7").Value = 163

$ws = $wb.Worksheets.Item(75)
$ws.Range("L4").Value = 29
$ws.Range("L7").Value = 207

$ws = $wb.Worksheets.Item(76)
$ws.Range("L4").Value = 5
$ws.Range("L7").Value = 52

$ws = $wb.Worksheets.Item(78)
$ws.Range("L3").Value = 37
$ws.Range("L7").Value = 105
